$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.70643966666667
$ws.Range("H2").Value = 53.119319
$ws.Range("I2").Value = 0.4380235920947999
$ws.Range("J2").Value = 0.4380235920947999
$ws.Range("M2").Value = 8.165540666666667
$ws.Range("N2").Value = 24.496622
$ws.Range("O2").Value = 0.1715865889461355
$ws.Range("P2").Value = 0.1715865889461355
$ws.Range("Q2").Value = 144.5826531600464
$ws.Range("R2").Value = 1301.243878440418
$ws.Range("S2").Value = 0.07515897404548014
$ws.Range("T2").Value = 0.07515897404548015
$ws.Range("G3").Value = 17.70643966666667
$ws.Range("H3").Value = 53.119319
$ws.Range("I3").Value = 0.4380235920947999
$ws.Range("J3").Value = 0.4380235920947999
$ws.Range("O3").Value = 0.4625449807101323
$ws.Range("P3").Value = 0.4625449807101323
$ws.Range("Q3").Value = 389.7506263611731
$ws.Range("R3").Value = 3507.755637250557
$ws.Range("S3").Value = 0.202605613956072
$ws.Range("T3").Value = 0.2026056139560721
$ws.Range("G4").Value = 17.70643966666667
$ws.Range("H4").Value = 53.119319
$ws.Range("I4").Value = 0.4380235920947999
$ws.Range("J4").Value = 0.4380235920947999
$ws.Range("M4").Value = 13.51552533333333
$ws.Range("N4").Value = 40.546576
$ws.Range("O4").Value = 0.2840084918355372
$ws.Range("P4").Value = 0.2840084918355373
$ws.Range("Q4").Value = 239.3118338779716
$ws.Range("R4").Value = 2153.806504901744
$ws.Range("S4").Value = 0.1244024197792287
$ws.Range("T4").Value = 0.1244024197792287
$ws.Range("G5").Value = 17.70643966666667
$ws.Range("H5").Value = 53.119319
$ws.Range("I5").Value = 0.4380235920947999
$ws.Range("J5").Value = 0.4380235920947999
$ws.Range("M5").Value = 3.895588
$ws.Range("N5").Value = 11.686764
$ws.Range("O5").Value = 0.08185993850819488
$ws.Range("P5").Value = 0.0818599385081949
$ws.Range("Q5").Value = 68.97699388819068
$ws.Range("R5").Value = 620.7929449937161
$ws.Range("S5").Value = 0.03585658431401895
$ws.Range("T5").Value = 0.03585658431401896
$ws.Range("G6").Value = 1.617245333333334
$ws.Range("H6").Value = 4.851736000000001
$ws.Range("I6").Value = 0.04000756919748267
$ws.Range("J6").Value = 0.04000756919748267
$ws.Range("M6").Value = 8.165540666666667
$ws.Range("N6").Value = 24.496622
$ws.Range("O6").Value = 0.1715865889461355
$ws.Range("P6").Value = 0.1715865889461355
$ws.Range("Q6").Value = 13.20568253731022
$ws.Range("R6").Value = 118.851142835792
$ws.Range("S6").Value = 0.006864762330622529
$ws.Range("T6").Value = 0.006864762330622531
$ws.Range("G7").Value = 1.617245333333334
$ws.Range("H7").Value = 4.851736000000001
$ws.Range("I7").Value = 0.04000756919748267
$ws.Range("J7").Value = 0.04000756919748267
$ws.Range("O7").Value = 0.4625449807101323
$ws.Range("P7").Value = 0.4625449807101323
$ws.Range("Q7").Value = 35.59848244551201
$ws.Range("R7").Value = 320.3863420096081
$ws.Range("S7").Value = 0.0185053003227089
$ws.Range("T7").Value = 0.01850530032270891
$ws.Range("G8").Value = 1.617245333333334
$ws.Range("H8").Value = 4.851736000000001
$ws.Range("I8").Value = 0.04000756919748267
$ws.Range("J8").Value = 0.04000756919748267
$ws.Range("M8").Value = 13.51552533333333
$ws.Range("N8").Value = 40.546576
$ws.Range("O8").Value = 0.2840084918355372
$ws.Range("P8").Value = 0.2840084918355373
$ws.Range("Q8").Value = 21.85792027288178
$ws.Range("R8").Value = 196.721282455936
$ws.Range("S8").Value = 0.01136248938978295
$ws.Range("T8").Value = 0.01136248938978295
$ws.Range("G9").Value = 1.617245333333334
$ws.Range("H9").Value = 4.851736000000001
$ws.Range("I9").Value = 0.04000756919748267
$ws.Range("J9").Value = 0.04000756919748267
$ws.Range("M9").Value = 3.895588
$ws.Range("N9").Value = 11.686764
$ws.Range("O9").Value = 0.08185993850819488
$ws.Range("P9").Value = 0.0818599385081949
$ws.Range("Q9").Value = 6.300121513589334
$ws.Range("R9").Value = 56.70109362230401
$ws.Range("S9").Value = 0.003275017154368283
$ws.Range("T9").Value = 0.003275017154368284
$ws.Range("G10").Value = 21.099799
$ws.Range("H10").Value = 63.299397
$ws.Range("I10").Value = 0.5219688387077175
$ws.Range("J10").Value = 0.5219688387077175
$ws.Range("M10").Value = 8.165540666666667
$ws.Range("N10").Value = 24.496622
$ws.Range("O10").Value = 0.1715865889461355
$ws.Range("P10").Value = 0.1715865889461355
$ws.Range("Q10").Value = 172.2912667929927
$ws.Range("R10").Value = 1550.621401136934
$ws.Range("S10").Value = 0.0895628525700328
$ws.Range("T10").Value = 0.08956285257003282
$ws.Range("G11").Value = 21.099799
$ws.Range("H11").Value = 63.299397
$ws.Range("I11").Value = 0.5219688387077175
$ws.Range("J11").Value = 0.5219688387077175
$ws.Range("O11").Value = 0.4625449807101323
$ws.Range("P11").Value = 0.4625449807101323
$ws.Range("Q11").Value = 464.444576727999
$ws.Range("R11").Value = 4180.001190551991
$ws.Range("S11").Value = 0.2414340664313513
$ws.Range("T11").Value = 0.2414340664313514
$ws.Range("G12").Value = 21.099799
$ws.Range("H12").Value = 63.299397
$ws.Range("I12").Value = 0.5219688387077175
$ws.Range("J12").Value = 0.5219688387077175
$ws.Range("M12").Value = 13.51552533333333
$ws.Range("N12").Value = 40.546576
$ws.Range("O12").Value = 0.2840084918355372
$ws.Range("P12").Value = 0.2840084918355373
$ws.Range("Q12").Value = 285.1748679127414
$ws.Range("R12").Value = 2566.573811214672
$ws.Range("S12").Value = 0.1482435826665256
$ws.Range("T12").Value = 0.1482435826665257
$ws.Range("G13").Value = 21.099799
$ws.Range("H13").Value = 63.299397
$ws.Range("I13").Value = 0.5219688387077175
$ws.Range("J13").Value = 0.5219688387077175
$ws.Range("M13").Value = 3.895588
$ws.Range("N13").Value = 11.686764
$ws.Range("O13").Value = 0.08185993850819488
$ws.Range("P13").Value = 0.0818599385081949
$ws.Range("Q13").Value = 82.19612378681201
$ws.Range("R13").Value = 739.765114081308
$ws.Range("S13").Value = 0.04272833703980765
$ws.Range("T13").Value = 0.04272833703980766
